# Add two new indicator rows (VRE penetration, Energy System Cost) to the
# id parameter table on Sheet1, along with their descriptions.
#
# Shared strings must be introduced in the order: both labels first, then
# both descriptions, so that the new shared-string table entries line up
# with the target workbook (labels occupy indices 147-148, descriptions
# occupy indices 149-150).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Labels first
$ws.Range("B76").Value = "VRE penetration"
$ws.Range("B77").Value = "Energy System Cost"

# Then descriptions
$ws.Range("C76").Value = "Variable Renewable Energy (VRE) penetration is defined as the share of wind and solar generation in annual gross electricity production"
$ws.Range("C77").Value = "Total energy system costs including balancing, profile, and grid costs, associated with integrating VRE at different penetration levels. [€/MWh]"

# Row ids
$ws.Range("A76").Value = 75
$ws.Range("A77").Value = 76

# Leave the selection where the author left it after entering the new data
$ws.Range("C79").Select() | Out-Null
